$d = $word.ActiveDocument

# Locate the "Full-Stack Development and Data Engineering" paragraph under the
# Siege Analytics / PARTNER entry - this is the anchor after which the new
# bullet paragraphs are inserted.
$anchorRange = $d.Content
$found = $anchorRange.Find.Execute("Full-Stack Development and Data Engineering",
                                    $true, $false, $false, $false, $false,
                                    $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find anchor paragraph 'Full-Stack Development and Data Engineering'"
}

# Get the actual Paragraph object so we have a stable Range to anchor inserts on.
$anchorPara = $anchorRange.Paragraphs.First
$anchor = $anchorPara.Range

$newBullets = @(
    "• Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States",
    "• Built scalable web applications processing 50,000+ electoral boundaries with sub-200ms response times",
    "• Architected systems supporting 2,500+ concurrent users conducting redistricting analysis",
    "• Algorithm reduced mapping costs by 75%, saving campaigns and organizations `$5M+ and enabling smaller nonprofits to conduct redistricting analysis"
)

# Insert one empty paragraph per bullet, directly after the anchor paragraph.
# Each call inserts right after the anchor (which is unaffected by the
# insertions happening after it), so the new paragraphs land in order.
foreach ($bullet in $newBullets) {
    $anchor.InsertParagraphAfter()
}

# Fill each of the newly-created empty paragraphs with its bullet text.
$firstNewIndex = $anchorPara.Index + 1
for ($i = 0; $i -lt $newBullets.Length; $i++) {
    $p = $d.Paragraphs($firstNewIndex + $i)
    $p.Range.InsertBefore($newBullets[$i])
}

Write-Output ("Done. ParaCount=" + $d.Paragraphs.Count)
